$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3-51 (B and C columns)
$bVals1 = @(-0.96,-0.919,-0.879,-0.838,-0.798,-0.758,-0.717,-0.677,-0.636,-0.596,-0.556,-0.515,-0.475,-0.434,-0.394,-0.354,-0.313,-0.273,-0.232,-0.192,-0.152,-0.111,-0.07099999999999999,-0.03,0.01,0.051,0.091,0.131,0.172,0.212,0.253,0.293,0.333,0.374,0.414,0.455,0.495,0.535,0.576,0.616,0.657,0.697,0.737,0.778,0.8179999999999999,0.859,0.899,0.9389999999999999,0.98)
$cVals1 = @(-3.587,-3.156,-2.858,-2.492,-2.131,-1.905,-1.676,-1.298,-1.138,1.512,-3.64,-0.369,-0.33,-0.046,0.177,0.276,0.316,0.532,0.665,0.7,0.867,0.8120000000000001,0.909,0.913,0.971,1.085,1.066,1.031,1.073,1.099,1.094,1.073,1.002,1.034,1.018,0.907,-2.137,0.77,0.776,0.667,0.596,0.623,-3.35,0.533,0.343,0.223,0.226,0.201,0.042)
for ($i = 0; $i -lt $bVals1.Length; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 2).Value = $bVals1[$i]
    $ws.Cells.Item($r, 3).Value = $cVals1[$i]
}

# Add new rows 52-101 (A, B, and C columns)
$aVals2 = @(50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99)
$bVals2 = @(1.02,1.061,1.101,1.141,1.182,1.222,1.263,1.303,1.343,1.384,1.424,1.465,1.505,1.545,1.586,1.626,1.667,1.707,1.747,1.788,1.828,1.869,1.909,1.949,1.99,2.03,2.071,2.111,2.152,2.192,2.232,2.273,2.313,2.354,2.394,2.434,2.475,2.515,2.556,2.596,2.636,2.677,2.717,2.758,2.798,2.838,2.879,2.919,2.96,3)
$cVals2 = @(-0.029,-0.091,-0.169,-0.355,-0.371,-3.445,-0.534,-0.63,-0.699,-0.6850000000000001,-0.695,-0.896,-0.856,-0.872,-3.737,-1.067,-0.954,-1.09,-1.016,-1.055,-1.067,-1.161,-1.107,-0.971,-4.672,-0.871,-0.983,-0.912,-0.705,-0.723,-3.167,-0.522,-0.305,-0.209,-0.127,0.063,0.211,0.483,-2.431,0.916,1.164,1.4,3.765,1.857,2.209,-0.491,2.8,3.198,3.678,4.048)
for ($i = 0; $i -lt $aVals2.Length; $i++) {
    $r = 52 + $i
    $ws.Cells.Item($r, 1).Value = $aVals2[$i]
    $ws.Cells.Item($r, 2).Value = $bVals2[$i]
    $ws.Cells.Item($r, 3).Value = $cVals2[$i]
}

# Copy the style of A51 (bold, bordered, centered) onto the new A column cells
$ws.Range("A51").Copy()
$ws.Range("A52:A101").PasteSpecial(-4122)
$excel.CutCopyMode = $false
